$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old sale rows (3-9), keeping only the header and row 2,
# which is about to be rewritten with the newest sale record.
for ($r = 9; $r -ge 3; $r--) {
    $ws.Rows.Item($r).Delete()
}

# Rewrite row 2 with the latest sale record.
$ws.Cells.Item(2, 1).Value = "2026-02-03 22:05:46"
$ws.Cells.Item(2, 3).Value = "Tarjeta de debito"
$ws.Cells.Item(2, 4).Value = "admin"

# Monto ("123.0") looks numeric, so a direct .Value assignment would be
# auto-coerced to a number (and would mint a brand new cell style for the
# text-forced cell). Build it as a text formula first, then paste-special
# the value back over itself so it collapses to a literal string while
# keeping the row's existing style untouched.
$ws.Cells.Item(2, 2).Formula = "=""123.0"""
$ws.Cells.Item(2, 2).Copy()
$ws.Cells.Item(2, 2).PasteSpecial(-4163)

# Narrow the Monto / Metodo de pago columns slightly.
$ws.Columns.Item(2).ColumnWidth = 6.15
$ws.Columns.Item(3).ColumnWidth = 18.15
